$d = $word.ActiveDocument

# Locate the unique phrase "según varios campos" inside the intro
# paragraph and apply a yellow highlight to it. This splits the original
# run into three runs (before / highlighted / after), matching the
# target diff.
$range = $d.Content
$found = $range.Find.Execute("según varios campos", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if ($found) {
    # wdYellow = 7
    $range.Font.HighlightColorIndex = 7
}
